# Requirements.txt for installing packages
# Update the "Contributions" slide (slide 8) to reflect that Mohammad Majid
# and James Bradford collaborated on both User Authentication and
# Javascript/Ajax/jQuery work.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# 1) "Mohammad Majid " -> "Mohammad Majid, James Bradford " (User Authentication line)
$full = $tr.Text
$idx = $full.IndexOf("Mohammad Majid ")
$sub = $tr.Characters($idx + 1, "Mohammad Majid ".Length)
$sub.Text = "Mohammad Majid, James Bradford "

# 2) "/Ajax:  " -> "/Ajax/jQuery:  " (Javascript line label)
$full = $tr.Text
$idx = $full.IndexOf("/Ajax:  ")
$sub = $tr.Characters($idx + 1, "/Ajax:  ".Length)
$sub.Text = "/Ajax/jQuery:  "

# 3) "James Bradford" (the occurrence right after the Ajax/jQuery label) ->
#    "James Bradford, Mohammad Majid". Search starting from the Ajax/jQuery
#    label so we don't match the "James Bradford" text introduced in step 1.
$full = $tr.Text
$ajaxIdx = $full.IndexOf("/Ajax/jQuery:  ")
$idx = $full.IndexOf("James Bradford", $ajaxIdx)
$sub = $tr.Characters($idx + 1, "James Bradford".Length)
$sub.Text = "James Bradford, Mohammad Majid"
